$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "21.731.02"
Set-TextValue $ws.Range("E2") "  -1.42%  "

Set-TextValue $ws.Range("D3") "1.540.54"
Set-TextValue $ws.Range("E3") "  -0.94%  "

Set-TextValue $ws.Range("E4") "  +0.09%  "

Set-TextValue $ws.Range("D5") "1.000"
Set-TextValue $ws.Range("E5") "  -0.03%  "

Set-TextValue $ws.Range("D6") "290.08"
Set-TextValue $ws.Range("E6") "  +1.12%  "

Set-TextValue $ws.Range("D7") "0.3881"
Set-TextValue $ws.Range("E7") "  +3.18%  "

Set-TextValue $ws.Range("D8") "0.3188"
Set-TextValue $ws.Range("E8") "  -1.69%  "

Set-TextValue $ws.Range("D9") "43.13"
Set-TextValue $ws.Range("E9") "  +5.05%  "

Set-TextValue $ws.Range("D10") "0.07205"
Set-TextValue $ws.Range("E10") "  -1.44%  "

Set-TextValue $ws.Range("D11") "1.062"
Set-TextValue $ws.Range("E11") "  -5.80%  "

Set-TextValue $ws.Range("D12") "0.9999"
Set-TextValue $ws.Range("E12") "  -0.03%  "

Set-TextValue $ws.Range("D13") "5.656"
Set-TextValue $ws.Range("E13") "  -0.97%  "

Set-TextValue $ws.Range("D14") "18.69"
Set-TextValue $ws.Range("E14") "  -4.55%  "

Set-TextValue $ws.Range("D15") "6.634"
Set-TextValue $ws.Range("E15") "  -3.29%  "

Set-TextValue $ws.Range("D16") "1.540.13"
Set-TextValue $ws.Range("E16") "  -1.04%  "

Set-TextValue $ws.Range("E17") "  +2.22%  "

Set-TextValue $ws.Range("D18") "0.06592"
Set-TextValue $ws.Range("E18") "  -0.86%  "

Set-TextValue $ws.Range("D19") "83.27"
Set-TextValue $ws.Range("E19") "  -2.28%  "

Set-TextValue $ws.Range("D20") "0.9999"
Set-TextValue $ws.Range("E20") "  -0.01%  "

Set-TextValue $ws.Range("D21") "6.165"
Set-TextValue $ws.Range("E21") "  -4.22%  "

Set-TextValue $ws.Range("D22") "15.42"
Set-TextValue $ws.Range("E22") "  -3.63%  "

Set-TextValue $ws.Range("D23") "10.91"
Set-TextValue $ws.Range("E23") "  -5.74%  "

Set-TextValue $ws.Range("D24") "2.407"
Set-TextValue $ws.Range("E24") "  +7.06%  "

Set-TextValue $ws.Range("D25") "21.736.96"
Set-TextValue $ws.Range("E25") "  -1.46%  "

Set-TextValue $ws.Range("D26") "2.378"
Set-TextValue $ws.Range("E26") "  -6.20%  "

Set-TextValue $ws.Range("D27") "146.56"
Set-TextValue $ws.Range("E27") "  -2.32%  "

Set-TextValue $ws.Range("D28") "18.43"
Set-TextValue $ws.Range("E28") "  -2.36%  "

Set-TextValue $ws.Range("D29") "4.837"
Set-TextValue $ws.Range("E29") "  +0.14%  "

Set-TextValue $ws.Range("D30") "1.716.31"
Set-TextValue $ws.Range("E30") "  -0.69%  "

Set-TextValue $ws.Range("D31") "117.73"
Set-TextValue $ws.Range("E31") "  -2.01%  "

Set-TextValue $ws.Range("D32") "0.9763"
Set-TextValue $ws.Range("E32") "  -13.11%  "

Set-TextValue $ws.Range("D33") "5.936"
Set-TextValue $ws.Range("E33") "  +0.05%  "

Set-TextValue $ws.Range("E34") "  +0.85%  "

Set-TextValue $ws.Range("D35") "8.845"
Set-TextValue $ws.Range("E35") "  -4.77%  "

Set-TextValue $ws.Range("D36") "0.06103"
Set-TextValue $ws.Range("E36") "  -0.72%  "

Set-TextValue $ws.Range("D37") "5.145"
Set-TextValue $ws.Range("E37") "  -1.68%  "

Set-TextValue $ws.Range("D38") "1.474"
Set-TextValue $ws.Range("E38") "  -13.39%  "

Set-TextValue $ws.Range("E39") "  -3.73%  "

Set-TextValue $ws.Range("D40") "0.2043"
Set-TextValue $ws.Range("E40") "  -3.45%  "

Set-TextValue $ws.Range("D41") "1.192"
Set-TextValue $ws.Range("E41") "  -2.21%  "

Set-TextValue $ws.Range("E42") "  +0.03%  "

Set-TextValue $ws.Range("D43") "10.70"
Set-TextValue $ws.Range("E43") "  -2.02%  "

Set-TextValue $ws.Range("E44") "  -3.26%  "

Set-TextValue $ws.Range("D45") "13.11"
Set-TextValue $ws.Range("E45") "  -3.62%  "

Set-TextValue $ws.Range("D46") "3.741"
Set-TextValue $ws.Range("E46") "  +0.30%  "

Set-TextValue $ws.Range("D47") "0.5515"
Set-TextValue $ws.Range("E47") "  -4.03%  "

Set-TextValue $ws.Range("D48") "1.162"
Set-TextValue $ws.Range("E48") "  +0.38%  "

Set-TextValue $ws.Range("D49") "116.90"
Set-TextValue $ws.Range("E49") "  -2.61%  "

Set-TextValue $ws.Range("D50") "1.878"
Set-TextValue $ws.Range("E50") "  -3.65%  "

Set-TextValue $ws.Range("E51") "  -3.04%  "
